# "updated the last scores" - refresh the SQUARNA-ver0.0C / SQUARNA-ver0.0B
# (columns L/M) benchmark numbers across the summary sheet, flattening the
# two F-score formula cells per dataset block into plain values (matching
# what Excel does when you type a literal over a formula cell), and tidy up
# a handful of stray empty/over-styled cells in column E / F left over from
# the SPOT-RNA column that has no data for the CoRToise datasets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TS1reduced block (rows 2-6) ---
$ws.Range("L3").Value = 0.604
$ws.Range("M3").Value = 0.631
$ws.Range("L4").Value = 0.613
$ws.Range("M4").Value = 0.633
$ws.Range("L5").Value = 0.736
$ws.Range("M5").Value = 0.76
$ws.Range("L6").Value = 0.538
$ws.Range("M6").Value = 0.554

# L2/M2 held formulas computing the F-score from L3:L4 / M3:M4; typing the
# refreshed score directly over them drops the formula and leaves a plain
# number, same as the recorded edit.
$ws.Range("L2").Value = 0.608
$ws.Range("M2").Value = 0.632

# --- CoRToise block (rows 12-16): drop stray empty, over-styled cells ---
$ws.Range("E12").Clear()
$ws.Range("E13").Clear()
$ws.Range("E14").Clear()
$ws.Range("E15").Clear()

# E5 / F15 carried a redundant duplicate style; restyle them to match the
# equivalent, already-used style (non-bold, centered) elsewhere on the sheet.
$ws.Range("A2").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("F15").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# E16 carried a redundant duplicate bold style; restyle it to match the
# equivalent bold/centered style used by the header row.
$ws.Range("A1").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# --- TS1reducedWC block (rows 17-21) ---
$ws.Range("L18").Value = 0.726
$ws.Range("M18").Value = 0.762
$ws.Range("L19").Value = 0.732
$ws.Range("M19").Value = 0.764
$ws.Range("L20").Value = 0.767
$ws.Range("M20").Value = 0.81
$ws.Range("L21").Value = 0.723
$ws.Range("M21").Value = 0.751

$ws.Range("L17").Value = 0.729
$ws.Range("M17").Value = 0.763

# Minor cosmetic drift recorded alongside the data refresh: selection moved
# to L10, and columns B/C widths nudged by a hundredth of a character.
$ws.Columns.Item(2).ColumnWidth = 11.04
$ws.Columns.Item(3).ColumnWidth = 8.21
$ws.Range("L10").Select()
